# Delete row 51 (the "Toy Story" / 「トイ・ストーリー」 post entry).
# This shifts all subsequent rows up by one, which matches the target
# workbook state (rows 52..214 become rows 51..213) and reduces the
# used range from A1:C214 to A1:C213.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(51).Delete()
